$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 14; A = "Globo"; B = "RJ TV 2"; C = "Social"; D = "2025-04-02T18:32"; E = "Negativo"; F = "Mutirão do CadÚnico em Campos. Ação de recadastramento é marcada por muita fila e reclamações. Entrevista com beneficiários. Reclamações de falta de banheiro e de água. Aplicativo não funciona. Eles querem voltar para o Cras. Imagens de fila grande. Ação começou hoje e segue até sexta. Entrevista com coordenadora do CadÚnico em Campos, Kamila Oliveira. O último mutirão aconteceu em fevereiro para beneficiários do Cartão Goitacá. 800 atendimentos por dia. Abertura dos portões 9h. *matéria* " },
    @{ Row = 15; A = "Globo"; B = "RJ TV 2"; C = "Social"; D = "2025-04-02T18:37"; E = "Negativo"; F = "Mutirão do Cadúnico em Campos. Moradores viram a madrugada buscando atendimento para atualização de dados. Repórter *ao vivo*. Reclamações continuam. Grande número de pessoas aguardando pelo atendimento de amanhã. Prefeitura por meio de nota informou que atendimento segue até sexta, mas não respondeu sobre reclamações feitas pelos beneficiários. " }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
